$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: add F35 and G35 values
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 0

# Updates to F/G columns across rows 271-401
$ws.Cells.Item(271, 6).Value = 45777
$ws.Cells.Item(272, 6).Value = 30778
$ws.Cells.Item(273, 6).Value = 31760
$ws.Cells.Item(273, 7).Value = 1665
$ws.Cells.Item(274, 6).Value = 28111
$ws.Cells.Item(274, 7).Value = 1279
$ws.Cells.Item(275, 6).Value = 30348
$ws.Cells.Item(276, 6).Value = 11345
$ws.Cells.Item(278, 6).Value = 30547
$ws.Cells.Item(278, 7).Value = 2102
$ws.Cells.Item(280, 6).Value = 34856
$ws.Cells.Item(280, 7).Value = 2316
$ws.Cells.Item(281, 6).Value = 46079
$ws.Cells.Item(281, 7).Value = 3164
$ws.Cells.Item(286, 6).Value = 55194
$ws.Cells.Item(286, 7).Value = 4282
$ws.Cells.Item(287, 6).Value = 58880
$ws.Cells.Item(287, 7).Value = 3719
$ws.Cells.Item(288, 6).Value = 59280
$ws.Cells.Item(288, 7).Value = 3975
$ws.Cells.Item(289, 6).Value = 62991
$ws.Cells.Item(289, 7).Value = 3592
$ws.Cells.Item(292, 6).Value = 82459
$ws.Cells.Item(292, 7).Value = 7275
$ws.Cells.Item(293, 6).Value = 82850
$ws.Cells.Item(293, 7).Value = 5776
$ws.Cells.Item(294, 6).Value = 93964
$ws.Cells.Item(294, 7).Value = 4949
$ws.Cells.Item(295, 6).Value = 17218
$ws.Cells.Item(295, 7).Value = 1036
$ws.Cells.Item(297, 6).Value = 2316
$ws.Cells.Item(299, 6).Value = 65746
$ws.Cells.Item(299, 7).Value = 6878
$ws.Cells.Item(300, 6).Value = 72562
$ws.Cells.Item(300, 7).Value = 6978
$ws.Cells.Item(301, 6).Value = 72198
$ws.Cells.Item(301, 7).Value = 5684
$ws.Cells.Item(305, 6).Value = 6756
$ws.Cells.Item(305, 7).Value = 526
$ws.Cells.Item(306, 6).Value = 149852
$ws.Cells.Item(306, 7).Value = 15228
$ws.Cells.Item(307, 6).Value = 151787
$ws.Cells.Item(307, 7).Value = 12798
$ws.Cells.Item(308, 6).Value = 30942
$ws.Cells.Item(308, 7).Value = 2100
$ws.Cells.Item(309, 6).Value = 155824
$ws.Cells.Item(309, 7).Value = 11062
$ws.Cells.Item(310, 6).Value = 158454
$ws.Cells.Item(311, 6).Value = 123012
$ws.Cells.Item(311, 7).Value = 3856
$ws.Cells.Item(312, 6).Value = 56270
$ws.Cells.Item(312, 7).Value = 1852
$ws.Cells.Item(313, 6).Value = 151191
$ws.Cells.Item(313, 7).Value = 6912
$ws.Cells.Item(314, 6).Value = 128738
$ws.Cells.Item(314, 7).Value = 6298
$ws.Cells.Item(315, 6).Value = 112740
$ws.Cells.Item(315, 7).Value = 5256
$ws.Cells.Item(316, 6).Value = 101504
$ws.Cells.Item(316, 7).Value = 4598
$ws.Cells.Item(317, 6).Value = 127479
$ws.Cells.Item(317, 7).Value = 4346
$ws.Cells.Item(318, 6).Value = 97928
$ws.Cells.Item(318, 7).Value = 2270
$ws.Cells.Item(320, 6).Value = 143192
$ws.Cells.Item(320, 7).Value = 6610
$ws.Cells.Item(321, 6).Value = 178713
$ws.Cells.Item(321, 7).Value = 5310
$ws.Cells.Item(322, 6).Value = 219365
$ws.Cells.Item(322, 7).Value = 4676
$ws.Cells.Item(323, 6).Value = 433677
$ws.Cells.Item(323, 7).Value = 6222
$ws.Cells.Item(324, 6).Value = 481621
$ws.Cells.Item(324, 7).Value = 5609
$ws.Cells.Item(325, 6).Value = 1531901
$ws.Cells.Item(325, 7).Value = 12918
$ws.Cells.Item(326, 6).Value = 839112
$ws.Cells.Item(326, 7).Value = 7670
$ws.Cells.Item(327, 6).Value = 450327
$ws.Cells.Item(327, 7).Value = 5450
$ws.Cells.Item(328, 6).Value = 361700
$ws.Cells.Item(328, 7).Value = 5342
$ws.Cells.Item(329, 6).Value = 165998
$ws.Cells.Item(329, 7).Value = 3518
$ws.Cells.Item(330, 6).Value = 145104
$ws.Cells.Item(330, 7).Value = 4168
$ws.Cells.Item(331, 6).Value = 309687
$ws.Cells.Item(331, 7).Value = 5409
$ws.Cells.Item(334, 6).Value = 196775
$ws.Cells.Item(334, 7).Value = 3504
$ws.Cells.Item(335, 6).Value = 130946
$ws.Cells.Item(335, 7).Value = 3001
$ws.Cells.Item(336, 6).Value = 102629
$ws.Cells.Item(336, 7).Value = 3344
$ws.Cells.Item(337, 6).Value = 103566
$ws.Cells.Item(337, 7).Value = 2891
$ws.Cells.Item(338, 6).Value = 228096
$ws.Cells.Item(338, 7).Value = 3201
$ws.Cells.Item(339, 6).Value = 661293
$ws.Cells.Item(339, 7).Value = 5494
$ws.Cells.Item(341, 6).Value = 291787
$ws.Cells.Item(341, 7).Value = 3667
$ws.Cells.Item(342, 6).Value = 178740
$ws.Cells.Item(342, 7).Value = 3039
$ws.Cells.Item(343, 6).Value = 133290
$ws.Cells.Item(343, 7).Value = 2973
$ws.Cells.Item(344, 6).Value = 135487
$ws.Cells.Item(344, 7).Value = 2487
$ws.Cells.Item(345, 6).Value = 292045
$ws.Cells.Item(345, 7).Value = 3328
$ws.Cells.Item(346, 6).Value = 675750
$ws.Cells.Item(346, 7).Value = 4835
$ws.Cells.Item(348, 6).Value = 232375
$ws.Cells.Item(348, 7).Value = 3238
$ws.Cells.Item(349, 6).Value = 158964
$ws.Cells.Item(349, 7).Value = 2752
$ws.Cells.Item(350, 6).Value = 127100
$ws.Cells.Item(350, 7).Value = 2787
$ws.Cells.Item(351, 6).Value = 150708
$ws.Cells.Item(351, 7).Value = 2831
$ws.Cells.Item(352, 6).Value = 307527
$ws.Cells.Item(352, 7).Value = 3546
$ws.Cells.Item(353, 6).Value = 725286
$ws.Cells.Item(353, 7).Value = 5298
$ws.Cells.Item(355, 6).Value = 221776
$ws.Cells.Item(355, 7).Value = 3444
$ws.Cells.Item(356, 6).Value = 159905
$ws.Cells.Item(356, 7).Value = 2876
$ws.Cells.Item(357, 6).Value = 138345
$ws.Cells.Item(357, 7).Value = 3026
$ws.Cells.Item(358, 6).Value = 158169
$ws.Cells.Item(358, 7).Value = 2605
$ws.Cells.Item(359, 6).Value = 321445
$ws.Cells.Item(359, 7).Value = 3348
$ws.Cells.Item(360, 6).Value = 750675
$ws.Cells.Item(360, 7).Value = 5137
$ws.Cells.Item(361, 6).Value = 332325
$ws.Cells.Item(361, 7).Value = 2619
$ws.Cells.Item(362, 6).Value = 228550
$ws.Cells.Item(362, 7).Value = 3181
$ws.Cells.Item(363, 6).Value = 187879
$ws.Cells.Item(363, 7).Value = 2761
$ws.Cells.Item(364, 6).Value = 167927
$ws.Cells.Item(364, 7).Value = 2473
$ws.Cells.Item(365, 6).Value = 184066
$ws.Cells.Item(365, 7).Value = 2393
$ws.Cells.Item(366, 6).Value = 339971
$ws.Cells.Item(366, 7).Value = 2852
$ws.Cells.Item(367, 6).Value = 764935
$ws.Cells.Item(367, 7).Value = 3914
$ws.Cells.Item(368, 6).Value = 345840
$ws.Cells.Item(368, 7).Value = 2295
$ws.Cells.Item(369, 6).Value = 233650
$ws.Cells.Item(369, 7).Value = 2595
$ws.Cells.Item(370, 6).Value = 182181
$ws.Cells.Item(370, 7).Value = 2039
$ws.Cells.Item(371, 6).Value = 159498
$ws.Cells.Item(371, 7).Value = 1956
$ws.Cells.Item(372, 6).Value = 179145
$ws.Cells.Item(372, 7).Value = 1862
$ws.Cells.Item(373, 6).Value = 349225
$ws.Cells.Item(373, 7).Value = 2375
$ws.Cells.Item(374, 6).Value = 771911
$ws.Cells.Item(375, 6).Value = 351286
$ws.Cells.Item(376, 6).Value = 220721
$ws.Cells.Item(376, 7).Value = 2222
$ws.Cells.Item(377, 6).Value = 176839
$ws.Cells.Item(377, 7).Value = 1821
$ws.Cells.Item(378, 6).Value = 157361
$ws.Cells.Item(378, 7).Value = 1546
$ws.Cells.Item(379, 6).Value = 179023
$ws.Cells.Item(380, 6).Value = 344057
$ws.Cells.Item(380, 7).Value = 2011
$ws.Cells.Item(381, 6).Value = 743851
$ws.Cells.Item(381, 7).Value = 2684
$ws.Cells.Item(383, 6).Value = 220790
$ws.Cells.Item(383, 7).Value = 1760
$ws.Cells.Item(384, 6).Value = 171986
$ws.Cells.Item(385, 6).Value = 150852
$ws.Cells.Item(385, 7).Value = 1404
$ws.Cells.Item(386, 6).Value = 182519
$ws.Cells.Item(386, 7).Value = 1361
$ws.Cells.Item(387, 6).Value = 351211
$ws.Cells.Item(388, 6).Value = 728039
$ws.Cells.Item(388, 7).Value = 2194
$ws.Cells.Item(390, 6).Value = 220499
$ws.Cells.Item(390, 7).Value = 1515
$ws.Cells.Item(391, 6).Value = 176573
$ws.Cells.Item(391, 7).Value = 1208
$ws.Cells.Item(392, 6).Value = 217506
$ws.Cells.Item(392, 7).Value = 1199
$ws.Cells.Item(393, 6).Value = 296915
$ws.Cells.Item(393, 7).Value = 1185
$ws.Cells.Item(394, 6).Value = 161621
$ws.Cells.Item(394, 7).Value = 614
$ws.Cells.Item(395, 6).Value = 734978
$ws.Cells.Item(395, 7).Value = 1912
$ws.Cells.Item(397, 6).Value = 106306
$ws.Cells.Item(397, 7).Value = 628
$ws.Cells.Item(398, 6).Value = 289297
$ws.Cells.Item(398, 7).Value = 1439
$ws.Cells.Item(399, 6).Value = 193810
$ws.Cells.Item(399, 7).Value = 955
$ws.Cells.Item(400, 6).Value = 142764
$ws.Cells.Item(400, 7).Value = 735
$ws.Cells.Item(401, 6).Value = 251986
$ws.Cells.Item(401, 7).Value = 883

# New row 402
$ws.Cells.Item(402, 1).Value = 44296
$ws.Cells.Item(402, 2).Value = 371062
$ws.Cells.Item(402, 3).Value = 7470
$ws.Cells.Item(402, 4).Value = 589
$ws.Cells.Item(402, 5).Value = 10565
$ws.Cells.Item(402, 6).Value = 610662
$ws.Cells.Item(402, 7).Value = 1162
